$d = $word.ActiveDocument

# 1) Fix the dash count in the "Attaching packages" verbatim line (add one dash).
$d.Content.Find.Execute(
    "## ── Attaching packages ──────────────────────────────────────────────────────────────────────────────────────── tidyverse 1.3.0 ──",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## ── Attaching packages ───────────────────────────────────────────────────────────────────────────────────────── tidyverse 1.3.0 ──",
    2) | Out-Null

# 2) Fix the dash count in the "Conflicts" verbatim line (add one dash).
$d.Content.Find.Execute(
    "## ── Conflicts ─────────────────────────────────────────────────────────────────────────────────────────── tidyverse_conflicts() ──",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## ── Conflicts ──────────────────────────────────────────────────────────────────────────────────────────── tidyverse_conflicts() ──",
    2) | Out-Null

# 3) Expand the "XX%" sentence in the lump/split paragraph with the replicated
#    (100x) calculation results. First swap in the full plain-text version,
#    then go back and bold each of the inserted numeric values so they end up
#    in their own runs, matching how Word splits formatting.
$oldSentence = "% were observed for the V4, V3-V4, and V4-V5 regions, respectively. However, at these thresholds, multiple species could be represented by the same OTU. At the highest level of resolution, XX% of the species shared a 16S rRNA gene sequence variant with another species. Given the risk of splitting a genome into multiple OTUs is more biologically problematic than lumping species together, larger thresholds are advisable."

$newSentence = "% were observed for the V4, V3-V4, and V4-V5 regions, respectively. However, at these thresholds, multiple species could be represented by the same OTU. At the highest level of resolution, 3.6% of the species shared a 16S rRNA gene sequence variant with another species when considering full length sequences and 14.9, 10.2, and 12.0% when considering the V4, V3-V4, and V4-V5 regions, respectively. At the commonly used 3% threshold, 25.2% of the species shared an OTU when considering full length sequences and 33.0, 29.4, and 32.2% when considering the V4, V3-V4, and V4-V5 regions, respectively. Given the risk of splitting a genome into multiple OTUs is more biologically problematic than lumping species together, larger thresholds are advisable."

$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null

$boldNumbers = @("3.6", "14.9", "10.2", "12.0", "25.2", "33.0", "29.4", "32.2")
foreach ($num in $boldNumbers) {
    $rng = $d.Content
    $found = $rng.Find.Execute($num, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Bold = 1

        # The number is always preceded by a single space in the source text.
        # Word keeps that space in its own (non-bold) run, separate from the
        # plain-text run before it -- nudge the engine into splitting it off
        # the same way by toggling a character property on/off.
        $spaceBefore = $d.Range($rng.Start - 1, $rng.Start)
        if ($spaceBefore.Text -eq " ") {
            $spaceBefore.Bold = 1
            $spaceBefore.Bold = 0
        }
    }
}
